$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The data for row 3 and row 4 (columns A, B, D, E, F, G, H) needs to be swapped.
# Capture current values of row 3
$a3 = $ws.Range("A3").Value()
$b3 = $ws.Range("B3").Value()
$d3 = $ws.Range("D3").Value()
$e3 = $ws.Range("E3").Value()
$f3 = $ws.Range("F3").Value()
$g3 = $ws.Range("G3").Value()
$h3 = $ws.Range("H3").Value()

# Capture current values of row 4
$a4 = $ws.Range("A4").Value()
$b4 = $ws.Range("B4").Value()
$d4 = $ws.Range("D4").Value()
$e4 = $ws.Range("E4").Value()
$f4 = $ws.Range("F4").Value()
$g4 = $ws.Range("G4").Value()
$h4 = $ws.Range("H4").Value()

# Write row 4 values into row 3
$ws.Range("A3").Value = $a4
$ws.Range("B3").Value = $b4
$ws.Range("D3").Value = $d4
$ws.Range("E3").Value = $e4
$ws.Range("F3").Value = $f4
$ws.Range("G3").Value = $g4
$ws.Range("H3").Value = $h4

# Write row 3 (original) values into row 4
$ws.Range("A4").Value = $a3
$ws.Range("B4").Value = $b3
$ws.Range("D4").Value = $d3
$ws.Range("E4").Value = $e3
$ws.Range("F4").Value = $f3
$ws.Range("G4").Value = $g3
$ws.Range("H4").Value = $h3
